$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: find the 4 plain paragraphs that make up the
#   public static void main(String[] args)
#   {
#       System.out.println("Hello World!");
#   }
# example block (inside the "Examples:" / "Hello World" section) and
# replace them with a VS-Code-style syntax-highlighted version.
# ---------------------------------------------------------------------
$cr = [char]13
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd($cr)
    if ($ptext -eq "public static void main(String[] args)") {
        $startIdx = $i
    }
    if ($startIdx -ge 0 -and $i -ge $startIdx -and $ptext -eq "}") {
        $endIdx = $i
        break
    }
}

if ($startIdx -lt 0 -or $endIdx -lt 0) {
    throw "Could not locate the Hello World example code block"
}

$blockStart = $d.Paragraphs.Item($startIdx).Range.Start
$blockEnd = $d.Paragraphs.Item($endIdx).Range.End
$blockRange = $d.Range($blockStart, $blockEnd)
$blockRange.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/><w:spacing w:line="270" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="569CD6"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>public</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="569CD6"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>static</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="4EC9B0"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>void</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="DCDCAA"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>main</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="4EC9B0"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>String</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">[] </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="9CDCFE"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>) {</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/><w:spacing w:line="270" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="9CDCFE"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>System</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="9CDCFE"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>out</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="DCDCAA"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>println</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="CE9178"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>&quot;Hello World&quot;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/><w:spacing w:line="270" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/><w:color w:val="D4D4D4"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">    }</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# ---------------------------------------------------------------------
# Step 2: move the hidden _GoBack bookmark from the "Next Assignment
# Link:" paragraph onto the start of the "In this example, the line"
# paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd($cr)
    if ($ptext -eq "In this example, the line") {
        $target = $d.Paragraphs.Item($i).Range
        $newRange = $d.Range($target.Start, $target.Start)
        $d.Bookmarks.Add("_GoBack", $newRange)
        break
    }
}
